$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws "D2" "91.745.63"
Set-TextValue $ws "E2" "  +0.54%  "
Set-TextValue $ws "D3" "3.090.79"
Set-TextValue $ws "E3" "  -1.99%  "
Set-TextValue $ws "E4" "  +0.02%  "
Set-TextValue $ws "D5" "235.33"
Set-TextValue $ws "E5" "  -1.77%  "
Set-TextValue $ws "D6" "611.76"
Set-TextValue $ws "E6" "  -1.26%  "
Set-TextValue $ws "E7" "  -3.40%  "
Set-TextValue $ws "E8" "  +2.40%  "
Set-TextValue $ws "D9" "1.00"
Set-TextValue $ws "E9" "  -0.02%  "
Set-TextValue $ws "D10" "3.085.85"
Set-TextValue $ws "E10" "  -2.07%  "
Set-TextValue $ws "D11" "0.729"
Set-TextValue $ws "E11" "  -1.88%  "
Set-TextValue $ws "E12" "  -1.65%  "
Set-TextValue $ws "E13" "  -0.95%  "
Set-TextValue $ws "D14" "91.924.73"
Set-TextValue $ws "E14" "  +0.96%  "
Set-TextValue $ws "D15" "33.82"
Set-TextValue $ws "E15" "  -3.91%  "
Set-TextValue $ws "D16" "5.40"
Set-TextValue $ws "E16" "  -3.46%  "
Set-TextValue $ws "D17" "3.683.47"
Set-TextValue $ws "E17" "  -1.67%  "
Set-TextValue $ws "D18" "3.093.81"
Set-TextValue $ws "E18" "  -1.85%  "
Set-TextValue $ws "E19" "  -2.47%  "
Set-TextValue $ws "D20" "14.51"
Set-TextValue $ws "E20" "  -5.06%  "
Set-TextValue $ws "D21" "5.74"
Set-TextValue $ws "E21" "  -5.07%  "
Set-TextValue $ws "B22" "BitcoinCash"
Set-TextValue $ws "C22" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws "D22" "442.26"
Set-TextValue $ws "E22" "  -2.34%  "
Set-TextValue $ws "B23" "Uniswap"
Set-TextValue $ws "C23" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws "D23" "9.23"
Set-TextValue $ws "E23" "  +0.35%  "
Set-TextValue $ws "E24" "  -5.57%  "
Set-TextValue $ws "D25" "5.71"
Set-TextValue $ws "E25" "  -4.96%  "
Set-TextValue $ws "D26" "85.72"
Set-TextValue $ws "E26" "  -3.66%  "
Set-TextValue $ws "D27" "11.54"
Set-TextValue $ws "E27" "  -4.31%  "
Set-TextValue $ws "D28" "3.259.07"
Set-TextValue $ws "E28" "  -1.98%  "
Set-TextValue $ws "E29" "  +0.01%  "
Set-TextValue $ws "E30" "  -1.49%  "
Set-TextValue $ws "D31" "0.231"
Set-TextValue $ws "E31" "  -2.36%  "
Set-TextValue $ws "E32" "  -2.65%  "
Set-TextValue $ws "E33" "  -3.74%  "
Set-TextValue $ws "B34" "RenderToken"
Set-TextValue $ws "C34" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue $ws "D34" "7.82"
Set-TextValue $ws "E34" "  +2.06%  "
Set-TextValue $ws "B35" "Kaspa"
Set-TextValue $ws "C35" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D35" "0.157"
Set-TextValue $ws "E35" "  -8.62%  "
Set-TextValue $ws "B36" "EthereumClassic"
Set-TextValue $ws "C36" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D36" "25.82"
Set-TextValue $ws "E36" "  -2.67%  "
Set-TextValue $ws "B37" "MantraDAO"
Set-TextValue $ws "C37" "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue $ws "D37" "3.90"
Set-TextValue $ws "E37" "  +0.72%  "
Set-TextValue $ws "B38" "PancakeSwap"
Set-TextValue $ws "C38" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws "D38" "1.88"
Set-TextValue $ws "E38" "  -4.15%  "
Set-TextValue $ws "B39" "Bittensor"
Set-TextValue $ws "C39" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws "D39" "479.23"
Set-TextValue $ws "E39" "  -6.57%  "
Set-TextValue $ws "B40" "WhiteBITCoin"
Set-TextValue $ws "C40" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws "D40" "23.85"
Set-TextValue $ws "E40" "  +7.76%  "
Set-TextValue $ws "B41" "Fetch.AI"
Set-TextValue $ws "C41" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws "D41" "1.28"
Set-TextValue $ws "E41" "  -5.19%  "
Set-TextValue $ws "B42" "PolygonEcosystemToken"
Set-TextValue $ws "C42" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws "D42" "0.427"
Set-TextValue $ws "E42" "  -5.10%  "
Set-TextValue $ws "B43" "Binance-PegBSC-USD"
Set-TextValue $ws "C43" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws "D43" "0.761"
Set-TextValue $ws "E43" "  -23.92%  "
Set-TextValue $ws "D44" "3.26"
Set-TextValue $ws "E44" "  -6.50%  "
Set-TextValue $ws "D46" "164.00"
Set-TextValue $ws "E46" "  +4.79%  "
Set-TextValue $ws "E47" "  -4.34%  "
Set-TextValue $ws "D48" "0.682"
Set-TextValue $ws "E48" "  -5.66%  "
Set-TextValue $ws "E49" "  -0.47%  "
Set-TextValue $ws "E50" "  +1.06%  "
Set-TextValue $ws "D51" "43.83"
Set-TextValue $ws "E51" "  -0.64%  "
